$d = $word.ActiveDocument

# Original sentence fragment (single run):
#   ", dentre as quais, serão utilizadas duas: Persona e Business Model "
# becomes:
#   ", dentre as quais, serão usadas três: Mapa da Empatia, Persona e Business Model "
# i.e. "utilizadas" -> "usadas", "duas" -> "três" and "Mapa da Empatia, " is
# inserted right before "Persona e Business Model ".

$old = ", dentre as quais, serão utilizadas duas: Persona e Business Model "
$new = ", dentre as quais, serão usadas três: Mapa da Empatia, Persona e Business Model "

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "replaced=$found"
